$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2 (shifts existing data rows 2-22 down to 3-23)
$ws.Rows.Item(2).Insert()

# Newly inserted rows inherit the formatting of the row above (the header);
# clear it so the new data row matches the plain, unstyled data rows below it.
$ws.Rows.Item(2).ClearFormats()

# Populate the newly inserted row with the new data point
$ws.Cells.Item(2, 1).Value = 0.003250675749898542
$ws.Cells.Item(2, 2).Value = 0.03200497691120401
$ws.Cells.Item(2, 3).Value = -0.08288132186446855

# Remove the old last two rows of data, which are now rows 22 and 23
$ws.Rows.Item(22).Resize(2).Delete()
